$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply additional AutoFilter criteria on top of the existing Academic Year
# filter: Year Level = 13, Qualification = "University Entrance".
# Field numbers are 1-based column offsets within the filter range (A=1, B=2, ... D=4).
$ws.Range("A1").AutoFilter(1, @("13"))
$ws.Range("A1").AutoFilter(2, @("University Entrance"))
$ws.Range("A1").AutoFilter(4, @("2018"))

# Move the active selection to H1
$ws.Range("H1").Select()
